# Update "想去人数" (interest count) figures in the F column
# on the "展览" and "全部类型" sheets, reflecting refreshed data
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value, for the "展览" sheet
$exhibitUpdates = @{
    2  = 21
    3  = 8329
    12 = 881
    13 = 3586
    14 = 249
    15 = 135
    16 = 778
    17 = 765
    19 = 484
    22 = 858
    23 = 1347
    24 = 397
    27 = 145
    28 = 323
    30 = 1014
    32 = 509
    33 = 624
    36 = 63
    39 = 129
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F-column value, for the "全部类型" sheet
$allTypesUpdates = @{
    3  = 21
    4  = 8329
    13 = 881
    15 = 3586
    16 = 249
    17 = 135
    19 = 778
    20 = 765
    23 = 484
    27 = 858
    28 = 1347
    29 = 397
    32 = 145
    34 = 323
    36 = 1014
    38 = 509
    39 = 624
    42 = 63
    45 = 129
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
